# (V1.2) Edit Pdf Model Diagram.
# Renames several "AddressBook"-style model classes to their "Pdf" equivalents
# on the class diagram, drops the now-unused Address field/connector, and
# refreshes the cached "last updated" date field everywhere it is cached
# (slide master, every slide layout, and the notes master).

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($container, [string]$newText)
    foreach ($sh in $container.Shapes) {
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# --- 1. Refresh the cached datetimeFigureOut field text everywhere ---------
Set-DatePlaceholderText $p.SlideMaster "3/9/19"

$customLayouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $customLayouts.Count; $i++) {
    Set-DatePlaceholderText $customLayouts.Item($i) "3/9/19"
}

Set-DatePlaceholderText $p.NotesMaster "3/9/19"

# --- 2. Rename model classes on the class diagram ---------------------------
$s = $p.Slides.Item(1)

function Set-ShapeTextById {
    param($slide, [int]$id, [string]$newText)
    foreach ($sh in $slide.Shapes) {
        if ($sh.Id -eq $id) {
            $sh.TextFrame.TextRange.Text = $newText
            return $sh
        }
    }
    return $null
}

function Get-ShapeById {
    param($slide, [int]$id)
    foreach ($sh in $slide.Shapes) {
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

[void](Set-ShapeTextById $s 46 "VersionedPdfBook")
[void](Set-ShapeTextById $s 49 "UniquePdfList")
[void](Set-ShapeTextById $s 62 "Pdf")
[void](Set-ShapeTextById $s 80 "Size")
[void](Set-ShapeTextById $s 100 "ReadOnlyPdfBook")
[void](Set-ShapeTextById $s 52 "Label")
[void](Set-ShapeTextById $s 55 "PdfBook")

# --- 3. Fold the "Email"/"Address" pair into a single "Location" field -----
# "Email" (id 83) becomes "Location" and shifts down very slightly (its
# position lines up with where "Address" used to be centred); the old
# "Address" rectangle (id 85) and its connector (id 84) are removed, leaving
# the remaining elbow connector (id 86) re-terminated on the shorter box.
$emailShape = Set-ShapeTextById $s 83 "Location"
if ($emailShape) {
    $emailShape.Top = 3210503 / 12700
}

$connector83 = Get-ShapeById $s 84
if ($connector83) { $connector83.Delete() }

$addressShape = Get-ShapeById $s 85
if ($addressShape) { $addressShape.Delete() }

$connector85 = Get-ShapeById $s 86
if ($connector85) {
    $connector85.Height = 318504 / 12700
}
